$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the label for the wealth factor: weight revised from 6% to 5%
$ws.Range("B3").Value = "Wealth (5%)"

# Leave the cursor/selection on the cell that was edited next (matches
# author's final selection state recorded in the sheet view)
$ws.Activate()
$ws.Range("B5").Select()
